$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 7.723999999999999
$ws.Range("B6").Value = 7.033999999999999
$ws.Range("B7").Value = 5.493
$ws.Range("C7").Value = -13.218
$ws.Range("B8").Value = 6.462000000000001
$ws.Range("C11").Value = -12.767
$ws.Range("C12").Value = -11.364
$ws.Range("D12").Value = -7.348999999999999
$ws.Range("D13").Value = -8.145
$ws.Range("D14").Value = -7.781999999999999
$ws.Range("C15").Value = -13.22
$ws.Range("B16").Value = 5.525
$ws.Range("D16").Value = -8.613999999999999
$ws.Range("D19").Value = -7.773999999999999
$ws.Range("B20").Value = 9.266999999999999
$ws.Range("C20").Value = -12.165
$ws.Range("D20").Value = -7.962000000000001
$ws.Range("B21").Value = 9.205
$ws.Range("C21").Value = -12.12
$ws.Range("C22").Value = -12.925
$ws.Range("D22").Value = -7.764
$ws.Range("C23").Value = -12.731
$ws.Range("B28").Value = 6.271999999999999
$ws.Range("B29").Value = 5.255000000000001
$ws.Range("C29").Value = -11.358
$ws.Range("B30").Value = 5.782
$ws.Range("B32").Value = 6.431
$ws.Range("C34").Value = -12.715
$ws.Range("D36").Value = -7.720000000000001
$ws.Range("B40").Value = 9.203999999999999
$ws.Range("C42").Value = -11.999
$ws.Range("C43").Value = -13.75
$ws.Range("D43").Value = -8.441999999999998
$ws.Range("C44").Value = -13.339
$ws.Range("C45").Value = -13.339
$ws.Range("B46").Value = 5.545
$ws.Range("C46").Value = -14.006
$ws.Range("D46").Value = -8.562000000000001
$ws.Range("C50").Value = -13.121
$ws.Range("D50").Value = -8.518000000000001
$ws.Range("B51").Value = 5.059
$ws.Range("C51").Value = -12.093
$ws.Range("B52").Value = 5.82
$ws.Range("B57").Value = 6.159000000000001
$ws.Range("C57").Value = -14.252
$ws.Range("B59").Value = 6.099000000000001
$ws.Range("B62").Value = 6.358
$ws.Range("C65").Value = -12.45
$ws.Range("B66").Value = 4.961
$ws.Range("C66").Value = -10.897
$ws.Range("C67").Value = -11.3
$ws.Range("B73").Value = 7.001
$ws.Range("B74").Value = 9.132999999999999
$ws.Range("D76").Value = -7.731999999999999
$ws.Range("B77").Value = 5.902000000000001
$ws.Range("C79").Value = -12.991
$ws.Range("C84").Value = -13.643
$ws.Range("C87").Value = -13.691
$ws.Range("B92").Value = 5.761
$ws.Range("C92").Value = -10.898
$ws.Range("D95").Value = -7.753
$ws.Range("C97").Value = -12.551
$ws.Range("D97").Value = -8.549000000000001
$ws.Range("D99").Value = -7.724999999999999
$ws.Range("B100").Value = 6.618
